# Add a "year 2000" data point to the top of each region's block on the
# "BLS Data Series" sheet (Northeast, Midwest, South, West), then leave the
# "Regions" sheet as the active tab with B14 selected on the data sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BLS Data Series")

# --- Northeast: insert new row 2 (year 2000) ---
$ws.Rows.Item(2).Insert()
$ws.Range("A3:C3").Copy()
$ws.Range("A2:C2").PasteSpecial(-4122)
$ws.Cells.Item(2, 1).Value = 2000
$ws.Cells.Item(2, 2).Value = "Northeast"
$ws.Cells.Item(2, 3).Value = 179.4

# --- Midwest: insert new row 23 (year 2000) ---
$ws.Rows.Item(23).Insert()
$ws.Range("A24:C24").Copy()
$ws.Range("A23:C23").PasteSpecial(-4122)
$ws.Cells.Item(23, 1).Value = 2000
$ws.Cells.Item(23, 2).Value = "Midwest"
$ws.Cells.Item(23, 3).Value = 168.3

# --- South: insert new row 44 (year 2000) ---
$ws.Rows.Item(44).Insert()
$ws.Range("A45:C45").Copy()
$ws.Range("A44:C44").PasteSpecial(-4122)
$ws.Cells.Item(44, 1).Value = 2000
$ws.Cells.Item(44, 2).Value = "South"
$ws.Cells.Item(44, 3).Value = 167.2

# --- West: insert new row 65 (year 2000) ---
$ws.Rows.Item(65).Insert()
$ws.Range("A66:C66").Copy()
$ws.Range("A65:C65").PasteSpecial(-4122)
$ws.Cells.Item(65, 1).Value = 2000
$ws.Cells.Item(65, 2).Value = "West"
$ws.Cells.Item(65, 3).Value = 174.8

$ws.Application.CutCopyMode = $false

# Selection left on the data sheet
$ws.Range("B14").Select()

# Make "Regions" the active (visible) tab, as in the saved workbook
$ws2 = $wb.Worksheets.Item("Regions")
$ws2.Activate()
